{"js": "const replacements = [\n  [\"42\u00d717=\", \"22\u00d778=\"],\n  [\"38\u00d790=\", \"36\u00d747=\"],\n  [\"22\u00d730=\", \"44\u00d715=\"],\n  [\"35\u00d773=\", \"44\u00d770=\"],\n  [\"95\u00d746=\", \"29\u00d736=\"],\n  [\"19\u00d757=\", \"27\u00d738=\"],\n  [\"73\u00d716=\", \"65\u00d783=\"],\n  [\"24\u00d788=\", \"27\u00d732=\"],\n  [\"14\u00d767=\", \"29\u00d726=\"],\n  [\"84\u00d782=\", \"43\u00d794=\"],\n  [\"25\u00d787=\", \"55\u00d750=\"],\n  [\"49\u00d764=\", \"29\u00d794=\"],\n  [\"53\u00d780=\", \"68\u00d796=\"],\n  [\"83\u00d773=\", \"14\u00d789=\"],\n  [\"54\u00d798=\", \"95\u00d779=\"],\n  [\"25\u00d735=\", \"31\u00d721=\"],\n  [\"52\u00d732=\", \"85\u00d728=\"],\n  [\"91\u00d737=\", \"67\u00d794=\"],\n  [\"12\u00d717=\", \"61\u00d717=\"],\n  [\"87\u00d731=\", \"54\u00d799=\"],\n  [\"66\u00d748=\", \"21\u00d713=\"],\n  [\"96\u00d719=\", \"12\u00d756=\"],\n  [\"54\u00d717=\", \"90\u00d799=\"],\n  [\"71\u00d773=\", \"74\u00d720=\"],\n  [\"11\u00d728=\", \"37\u00d712=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @{Old=\"42\u00d717=\"; New=\"22\u00d778=\"},\n  @{Old=\"38\u00d790=\"; New=\"36\u00d747=\"},\n  @{Old=\"22\u00d730=\"; New=\"44\u00d715=\"},\n  @{Old=\"35\u00d773=\"; New=\"44\u00d770=\"},\n  @{Old=\"95\u00d746=\"; New=\"29\u00d736=\"},\n  @{Old=\"19\u00d757=\"; New=\"27\u00d738=\"},\n  @{Old=\"73\u00d716=\"; New=\"65\u00d783=\"},\n  @{Old=\"24\u00d788=\"; New=\"27\u00d732=\"},\n  @{Old=\"14\u00d767=\"; New=\"29\u00d726=\"},\n  @{Old=\"84\u00d782=\"; New=\"43\u00d794=\"},\n  @{Old=\"25\u00d787=\"; New=\"55\u00d750=\"},\n  @{Old=\"49\u00d764=\"; New=\"29\u00d794=\"},\n  @{Old=\"53\u00d780=\"; New=\"68\u00d796=\"},\n  @{Old=\"83\u00d773=\"; New=\"14\u00d789=\"},\n  @{Old=\"54\u00d798=\"; New=\"95\u00d779=\"},\n  @{Old=\"25\u00d735=\"; New=\"31\u00d721=\"},\n  @{Old=\"52\u00d732=\"; New=\"85\u00d728=\"},\n  @{Old=\"91\u00d737=\"; New=\"67\u00d794=\"},\n  @{Old=\"12\u00d717=\"; New=\"61\u00d717=\"},\n  @{Old=\"87\u00d731=\"; New=\"54\u00d799=\"},\n  @{Old=\"66\u00d748=\"; New=\"21\u00d713=\"},\n  @{Old=\"96\u00d719=\"; New=\"12\u00d756=\"},\n  @{Old=\"54\u00d717=\"; New=\"90\u00d799=\"},\n  @{Old=\"71\u00d773=\"; New=\"74\u00d720=\"},\n  @{Old=\"11\u00d728=\"; New=\"37\u00d712=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.Text = $r.Old\n    $range.Find.Replacement.Text = $r.New\n    $range.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
